# Generate Report for Handoff
# Update the "b.md" rows across the Overview / zh-cn / de-de sheets to
# reflect that b.md is now ready for handoff, with fresh handoff file
# names (b.<hash>.<locale>.xlf) and new handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row (A3 = "b.md").
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-22 10:33:36"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row (A3 = "b.md").
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-22 10:33:32"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Row -eq 3 -and $h.Range.Column -eq 4) {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row (A3 = "b.md").
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-22 10:33:36"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Row -eq 3 -and $h.Range.Column -eq 4) {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
